$d = $word.ActiveDocument

# --- Paragraph 6 (last): "spel 4: Boss fight" -> "Spel 4: " + "ontsnap het doolhof..." ---
$p6 = $d.Paragraphs(6)
$xml6 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Spel </w:t></w:r><w:r><w:t>4</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p><w:r><w:t>ontsnap het doolhof, ren door het doolhof met dezelfde knoppen als de map</w:t></w:r><w:r><w:t xml:space="preserve">, gebruik de het nieuwe gereedschap (knop 5 van links voor de hamer en knop 6 voor de poolstok) om door obstakels te komen</w:t></w:r><w:r><w:t xml:space="preserve"> en ontsnap eindelijk!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p6.Range.InsertXML($xml6)
# Remove the now-superfluous trailing empty paragraph left behind because this
# was the final paragraph in the document (its mark carries the body sectPr).
$secondLast = $d.Paragraphs($d.Paragraphs.Count - 1)
$d.Range($secondLast.Range.End - 1, $secondLast.Range.End).Delete()

# --- Paragraph 5: "spel 3: Jumping game on LCD screen" -> "spel 3: " + "Parcour!..." ---
$p5 = $d.Paragraphs(5)
$xml5 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">spel 3: </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Parcour</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>! Wissel tussen de bovenkant en onderkant van de rijen op het lcd scherm en bereik het einde om de poolstok te krijgen waarmee je over gaten kunt springen.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p5.Range.InsertXML($xml5)

# --- Paragraph 4: "spel 2: Whack a mole..." -> "spel 2: " + "Whack a mole...slopen." ---
$p4 = $d.Paragraphs(4)
$xml4 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">spel </w:t></w:r><w:r><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Whack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mole</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. Druk de verschillende lampjes zullen willekeurig gaan branden. Druk de knoppen onder de lampjes die beginnen te branden op tijd aan om punten te scoren. Elke fout haalt 2 punten weg. Bereik de 2</w:t></w:r><w:r><w:t>5</w:t></w:r><w:r><w:t xml:space="preserve"> punten en win. Het spel zal moeilijker worden hoe dichter je bij het einde komt</w:t></w:r><w:r><w:t xml:space="preserve">! </w:t></w:r><w:r><w:br/><w:t>Als je het spel voltooid krijg je de hamer die obstakels in het doolhof kan slopen.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p4.Range.InsertXML($xml4)

# --- Paragraph 3: "spel 1: Balanceren...seconden" -> "Navigeer..." + "spel 1: " + "Balanceren...wijst!" ---
$p3 = $d.Paragraphs(3)
$xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Navigeer rond de map om de spellen te vinden. Gebruik de eerste knop links of naar links te bewegen, de tweede knop om naar rechts te bewegen, de derde knop om omhoog te bewegen en de vierde knop om omlaag te bewegen. </w:t></w:r><w:r><w:t>Voltooi alle spellen om items te krijgen en via het doolhof te kunnen ontsnappen!</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">spel </w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p><w:r><w:t>Balanceren. Gebruik de potentiometer om het brandende lampje in het midden te houden voor 10 seconden</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">Als je het spel voltooid heb krijg je het </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>compas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> dat je richting het einde wijst!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p3.Range.InsertXML($xml3)
